$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "27.366.10"
$ws.Range("E2").Value = "  +1.64%  "
$ws.Range("D3").Value = "1.844.93"
$ws.Range("E3").Value = "  +1.41%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value = "1.014"
$ws.Range("E4").Value = "  +1.29%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "315.68"
$ws.Range("E5").Value = "  +1.97%  "
$ws.Range("E6").Value = "  +1.04%  "
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = "0.4738"
$ws.Range("E7").Value = "  +1.16%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "0.3703"
$ws.Range("E8").Value = "  +0.21%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.07461"
$ws.Range("E9").Value = "  +1.20%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.8870"
$ws.Range("E10").Value = "  +1.66%  "
$ws.Range("E11").Value = "  +0.12%  "
$ws.Range("D12").Value = "1.860.08"
$ws.Range("E12").Value = "  +2.47%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "0.07371"
$ws.Range("E13").Value = "  +4.23%  "
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "5.486"
$ws.Range("E14").Value = "  +2.17%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "93.53"
$ws.Range("E15").Value = "  +0.95%  "
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = "6.590"
$ws.Range("E17").Value = "  +1.21%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "0.000008852"
$ws.Range("E18").Value = "  +1.49%  "
$ws.Range("E19").Value = "  +0.98%  "
$ws.Range("E20").Value = "  +0.51%  "
$ws.Range("D21").Value = "27.394.46"
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = "5.340"
$ws.Range("E23").Value = "  +1.34%  "
$ws.Range("D24").Value = "2.066.13"
$ws.Range("E24").Value = "  +1.30%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "1.906"
$ws.Range("E25").Value = "  +0.22%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "152.94"
$ws.Range("E26").Value = "  +0.98%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "18.69"
$ws.Range("E27").Value = "  +1.57%  "
$ws.Range("E28").Value = "  +0.12%  "
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = "5.303"
$ws.Range("E29").Value = "  -0.51%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = "118.28"
$ws.Range("E30").Value = "  +2.12%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.08970"
$ws.Range("E31").Value = "  +0.33%  "
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "0.7645"
$ws.Range("E32").Value = "  -0.93%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = "1.180"
$ws.Range("E33").Value = "  +1.10%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = "4.575"
$ws.Range("E34").Value = "  +1.58%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "2.946"
$ws.Range("E35").Value = "  +1.51%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "1.011"
$ws.Range("E36").Value = "  +1.05%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "1.110"
$ws.Range("E37").Value = "  +2.00%  "
$ws.Range("E38").Value = "  +1.37%  "
$ws.Range("E39").Value = "  -0.03%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "3.008"
$ws.Range("E40").Value = "  +1.73%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "7.364"
$ws.Range("E41").Value = "  +0.68%  "
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = "2.418"
$ws.Range("E42").Value = "  +2.29%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = "0.5373"
$ws.Range("E43").Value = "  +0.42%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = "0.1669"
$ws.Range("E44").Value = "  -0.03%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "8.584"
$ws.Range("E45").Value = "  +1.63%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "0.4976"
$ws.Range("E46").Value = "  +0.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "10.58"
$ws.Range("E47").Value = "  +1.19%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "1.685"
$ws.Range("E49").Value = "  +0.83%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value = "104.55"
$ws.Range("E50").Value = "  +0.63%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = "0.06330"
$ws.Range("E51").Value = "  +0.70%  "
